# "Changes in method diagnosis"
#
# MEDICAL EXAMINATIONS table (B26:E31) gains a new "patient_id" column (F):
#   F27 = header "patient_id" (same look as the other yellow-filled header
#         cells in that row, i.e. E27's border/alignment + B27's yellow fill)
#   F28:F31 = the id_patient values for each exam row (same style as the
#         other data cells in that row, e.g. C28:C31)
#
# The SYMPTOMPS table (B34:J38) used to re-key on id_medExam (column D) -
# those link values are removed (column D is now blank, header text stays).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- MEDICAL EXAMINATIONS: new "patient_id" column header (F27) ---------
# Start from the same format as the neighbouring header cell (E27), then
# recolor it to match the rest of the yellow header row (B27).
$ws.Range("E27").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Interior.Color = 65535
$ws.Range("F27").Value = "patient_id"

# --- MEDICAL EXAMINATIONS: new "patient_id" values (F28:F31) -------------
$ws.Range("C28").Copy()
$ws.Range("F28:F31").PasteSpecial(-4122)
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("F30").Value = 2
$ws.Range("F31").Value = 6

$excel.CutCopyMode = $false

# --- SYMPTOMPS: clear the old id_medExam link values (D36:D38) -----------
$ws.Range("D35:D38").ClearContents()

# --- Window/selection state, matching where the author left the cursor ---
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
[void]$ws.Range("D35").Select()
